# Complete ramp up function
# Applies the commit's changes:
#  - Flight Mission Cycle: Light switch "No. of cycles" 1 -> 3; remove the Piano row
#  - Typing: RoM type sinosoidal -> triangle; Max_RoM 10 -> 0; Min_RoM -20 -> 0; Period 30 -> 15
#  - Light switch: Min_RoM -20 -> 10
#  - Active sheet moves from "Light switch" to "Typing"; per-sheet selections updated

$wb = $excel.ActiveWorkbook

# --- Flight Mission Cycle sheet ---
$ws = $wb.Worksheets.Item("Flight Mission Cycle")
[void]$ws.Select()
$ws.Range("B3").Value = 3
[void]$ws.Range("A4:B4").EntireRow.Delete()
[void]$ws.Range("E12").Select()

# --- Typing sheet ---
$wsTyping = $wb.Worksheets.Item("Typing")
[void]$wsTyping.Select()
$wsTyping.Range("B4").Value = "triangle"
$wsTyping.Range("C5").Value = 0
$wsTyping.Range("C6").Value = 0
$wsTyping.Range("C7").Value = 15

# --- Light switch sheet ---
$wsLightSwitch = $wb.Worksheets.Item("Light switch")
[void]$wsLightSwitch.Select()
$wsLightSwitch.Range("C6").Value = 10
[void]$wsLightSwitch.Range("B2").Select()

# Finish with "Typing" selected/active (activeTab moves from Light switch to
# Typing) and its last selection on H12, matching the saved view state.
[void]$wsTyping.Select()
[void]$wsTyping.Range("H12").Select()
